$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the title banner text (A1) - shared string changes from 05.05.2020 to 06.05.2020
$ws.Range("A1").Value = "Données COVID-19 Valais 06.05.2020"

# 2) Update the selection shown on the sheet view (A1:L1, matching the title merge range)
$ws.Range("A1:L1").Select()

# 3) Update existing rows 59-69 with revised figures (G column mostly, plus C69/D68)
$ws.Range("G59").Value = 57
$ws.Range("G60").Value = 60
$ws.Range("G61").Value = 61
$ws.Range("G62").Value = 52
$ws.Range("G63").Value = 49
$ws.Range("G64").Value = 50
$ws.Range("G65").Value = 51
$ws.Range("G66").Value = 49
$ws.Range("G67").Value = 49
$ws.Range("D68").Value = 1
$ws.Range("G68").Value = 49
$ws.Range("C69").Value = 7
$ws.Range("G69").Value = 44

# 4) Insert a fresh row at position 70 - this pushes the former "last" data row (70) down
#    to 71, and the trailing footnote row (71) down to 72, exactly matching the diff's
#    row-shift.
$ws.Rows("70").Insert()

# 5) Give the newly inserted row 70 the same formatting as row 69 (the row above it).
$ws.Range("A69:L69").Copy()
$ws.Range("A70:L70").PasteSpecial(-4122)
$ws.Range("A70:L70").Select()
$excel.CutCopyMode = $false

# 6) Populate the new row 70 (regular data row)
$ws.Range("A70").Value = 43956
$ws.Range("B70").Formula = "=B69+C70"
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = 6
$ws.Range("G70").Value = 41
$ws.Range("H70").Formula = "=G70+E70"
$ws.Range("I70").Formula = "=I69+J70"
$ws.Range("J70").Formula = "=K70+L70"

# K70/L70 inherit a text number format from the copied row; flip to General so the
# values land as real numbers (matching the rest of the K/L column), then restore
# the original "text" display format those columns use.
$kFmt = $ws.Range("K70").NumberFormat
$lFmt = $ws.Range("L70").NumberFormat
$ws.Range("K70").NumberFormat = "General"
$ws.Range("L70").NumberFormat = "General"
$ws.Range("K70").Value = 2
$ws.Range("L70").Value = 2
$ws.Range("K70").NumberFormat = $kFmt
$ws.Range("L70").NumberFormat = $lFmt

# 7) Row 71 already carries the former "final row" style (it was shifted down by the
#    insert above) - overwrite its values with the updated figures.
$ws.Range("A71").Value = 43957
$ws.Range("B71").ClearContents()
$ws.Range("C71").ClearContents()
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 6
$ws.Range("G71").Value = 41
$ws.Range("H71").Formula = "=G71+E71"
$ws.Range("I71").Formula = "=I70+J71"
$ws.Range("J71").Formula = "=K71+L71"
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
